{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Map of exact current paragraph text -> new paragraph text.\nconst replacements = new Map([\n  [\"New in this update (Render deployment setup)\", \"New in this update (Render npm ERESOLVE fix)\"],\n  [\"- Added Render Blueprint file: `render.yaml`.\", \"- Fixed frontend dependency conflict causing Render build failure:\"],\n  [\"  - Provisions PostgreSQL, backend service, and frontend static service.\", \"  - Changed `date-fns` from `^4.1.0` to `^3.6.0` in `frontend/package.json`.\"],\n  [\"- Added Render deploy runbook: `DEPLOY_RENDER.md`.\", \"  - This satisfies `react-day-picker@8.10.1` peer requirement (`^2.28.0 || ^3.0.0`).\"],\n  [\"- Render build config includes frontend install fallback:\", \"- Updated `render.yaml` frontend build command back to:\"],\n  [\"  - `npm install --legacy-peer-deps && npm run build`.\", \"  - `npm install && npm run build`\"],\n  [\"- Last pushed commit: 09217f9\", \"- Last pushed commit: 0369be1\"],\n  [\"- Current Render deployment setup is local and not pushed yet.\", \"- Current Render dependency fix is local and not pushed yet.\"],\n]);\n\nfor (const p of paragraphs.items) {\n  if (replacements.has(p.text)) {\n    p.insertText(replacements.get(p.text), \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-ExactParagraphText($doc, [string]$oldText, [string]$newText) {\n  foreach ($p in $doc.Paragraphs) {\n    $r = $p.Range\n    # Paragraph range text includes trailing paragraph mark; strip it for comparison.\n    $t = $r.Text\n    if ($t.Length -gt 0 -and ($t.EndsWith([char]13) -or $t.EndsWith([char]7))) {\n      $t = $t.Substring(0, $t.Length - 1)\n    }\n    if ($t -eq $oldText) {\n      $r.Text = $newText\n      return\n    }\n  }\n}\n\nReplace-ExactParagraphText $d \"New in this update (Render deployment setup)\" \"New in this update (Render npm ERESOLVE fix)\"\nReplace-ExactParagraphText $d \"- Added Render Blueprint file: ``render.yaml``.\" \"- Fixed frontend dependency conflict causing Render build failure:\"\nReplace-ExactParagraphText $d \"  - Provisions PostgreSQL, backend service, and frontend static service.\" \"  - Changed ``date-fns`` from ``^4.1.0`` to ``^3.6.0`` in ``frontend/package.json``.\"\nReplace-ExactParagraphText $d \"- Added Render deploy runbook: ``DEPLOY_RENDER.md``.\" \"  - This satisfies ``react-day-picker@8.10.1`` peer requirement (``^2.28.0 || ^3.0.0``).\"\nReplace-ExactParagraphText $d \"- Render build config includes frontend install fallback:\" \"- Updated ``render.yaml`` frontend build command back to:\"\nReplace-ExactParagraphText $d \"  - ``npm install --legacy-peer-deps && npm run build``.\" \"  - ``npm install && npm run build``\"\nReplace-ExactParagraphText $d \"- Last pushed commit: 09217f9\" \"- Last pushed commit: 0369be1\"\nReplace-ExactParagraphText $d \"- Current Render deployment setup is local and not pushed yet.\" \"- Current Render dependency fix is local and not pushed yet.\"\n"}
